{"js": "const replacements = [\n  [\"2024-08-30 Friday\", \"2024-08-31 Saturday\"],\n  [\"235\u00f77=\", \"709\u00f74=\"],\n  [\"361\u00f76=\", \"983\u00f79=\"],\n  [\"519\u00f77=\", \"221\u00f77=\"],\n  [\"421\u00f75=\", \"287\u00f75=\"],\n  [\"575\u00f77=\", \"532\u00f76=\"],\n  [\"854\u00f75=\", \"863\u00f74=\"],\n  [\"689\u00f74=\", \"792\u00f73=\"],\n  [\"684\u00f79=\", \"403\u00f79=\"],\n  [\"375\u00f77=\", \"593\u00f79=\"],\n  [\"179\u00f77=\", \"912\u00f79=\"],\n  [\"474\u00f73=\", \"317\u00f72=\"],\n  [\"258\u00f75=\", \"436\u00f76=\"],\n  [\"998\u00f78=\", \"147\u00f75=\"],\n  [\"170\u00f77=\", \"933\u00f77=\"],\n  [\"999\u00f72=\", \"226\u00f74=\"],\n  [\"246\u00f73=\", \"335\u00f79=\"],\n  [\"952\u00f76=\", \"449\u00f78=\"],\n  [\"822\u00f79=\", \"789\u00f79=\"],\n  [\"529\u00f78=\", \"597\u00f74=\"],\n  [\"986\u00f79=\", \"799\u00f72=\"],\n  [\"968\u00f73=\", \"509\u00f72=\"],\n  [\"871\u00f72=\", \"307\u00f74=\"],\n  [\"803\u00f79=\", \"572\u00f74=\"],\n  [\"433\u00f73=\", \"729\u00f77=\"],\n  [\"371\u00f76=\", \"666\u00f78=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-08-30 Friday', '2024-08-31 Saturday'),\n    @('235\u00f77=', '709\u00f74='),\n    @('361\u00f76=', '983\u00f79='),\n    @('519\u00f77=', '221\u00f77='),\n    @('421\u00f75=', '287\u00f75='),\n    @('575\u00f77=', '532\u00f76='),\n    @('854\u00f75=', '863\u00f74='),\n    @('689\u00f74=', '792\u00f73='),\n    @('684\u00f79=', '403\u00f79='),\n    @('375\u00f77=', '593\u00f79='),\n    @('179\u00f77=', '912\u00f79='),\n    @('474\u00f73=', '317\u00f72='),\n    @('258\u00f75=', '436\u00f76='),\n    @('998\u00f78=', '147\u00f75='),\n    @('170\u00f77=', '933\u00f77='),\n    @('999\u00f72=', '226\u00f74='),\n    @('246\u00f73=', '335\u00f79='),\n    @('952\u00f76=', '449\u00f78='),\n    @('822\u00f79=', '789\u00f79='),\n    @('529\u00f78=', '597\u00f74='),\n    @('986\u00f79=', '799\u00f72='),\n    @('968\u00f73=', '509\u00f72='),\n    @('871\u00f72=', '307\u00f74='),\n    @('803\u00f79=', '572\u00f74='),\n    @('433\u00f73=', '729\u00f77='),\n    @('371\u00f76=', '666\u00f78='),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}"}
